$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "湖南黄金"
$ws.Range("B2").Value = "白银有色"
$ws.Range("C2").Value = "白银有色"
$ws.Range("A3").Value = "白银有色"
$ws.Range("B3").Value = "红 宝 丽"
$ws.Range("C3").Value = "利欧股份"
$ws.Range("A4").Value = "红 宝 丽"
$ws.Range("B4").Value = "湖南黄金"
$ws.Range("C4").Value = "红宝丽"
$ws.Range("A5").Value = "利欧股份"
$ws.Range("B5").Value = "贵州茅台"
$ws.Range("C5").Value = "天奇股份"
$ws.Range("A6").Value = "新易盛"
$ws.Range("B6").Value = "蓝色光标"
$ws.Range("C6").Value = "航天发展"
$ws.Range("A7").Value = "蓝色光标"
$ws.Range("B7").Value = "利欧股份"
$ws.Range("C7").Value = "湖南黄金"
$ws.Range("A8").Value = "中际旭创"
$ws.Range("B8").Value = "东方财富"
$ws.Range("C8").Value = "蓝色光标"
$ws.Range("A9").Value = "湖南白银"
$ws.Range("B9").Value = "中国黄金"
$ws.Range("C9").Value = "紫金矿业"
$ws.Range("A10").Value = "紫金矿业"
$ws.Range("B10").Value = "湖南白银"
$ws.Range("C10").Value = "湖南白银"
$ws.Range("A11").Value = "天奇股份"
$ws.Range("B11").Value = "亨通光电"
$ws.Range("C11").Value = "新易盛"
$ws.Range("A12").Value = "中国黄金"
$ws.Range("B12").Value = "紫金矿业"
$ws.Range("C12").Value = "中际旭创"
$ws.Range("A13").Value = "航天发展"
$ws.Range("B13").Value = "农发种业"
$ws.Range("C13").Value = "农发种业"
$ws.Range("A14").Value = "农发种业"
$ws.Range("B14").Value = "新易盛"
$ws.Range("C14").Value = "天地在线"
$ws.Range("A15").Value = "亨通光电"
$ws.Range("B15").Value = "天奇股份"
$ws.Range("C15").Value = "山子高科"
$ws.Range("A16").Value = "贵州茅台"
$ws.Range("B16").Value = "太极实业"
$ws.Range("C16").Value = "中国黄金"
$ws.Range("A17").Value = "西部材料"
$ws.Range("B17").Value = "通鼎互联"
$ws.Range("C17").Value = "特变电工"
$ws.Range("A18").Value = "铜陵有色"
$ws.Range("B18").Value = "铜陵有色"
$ws.Range("C18").Value = "铜陵有色"
$ws.Range("A19").Value = "华天科技"
$ws.Range("B19").Value = "航天发展"
$ws.Range("C19").Value = "华天科技"
$ws.Range("A20").Value = "信维通信"
$ws.Range("B20").Value = "百川股份"
$ws.Range("C20").Value = "亨通光电"
$ws.Range("A21").Value = "天地在线"
$ws.Range("B21").Value = "华天科技"
$ws.Range("C21").Value = "巨力索具"
